# xlsx data to db
#
# Normalize the "Branch" column against the DB values:
#   D2 (Apoorv Aron):    IT  -> B.Tech IT
#   D3 (Vishesh Gupta):   CSE -> B.Tech CSE
#   D4 (Nidhi Rathore):   IT  -> MCA
#
# Vishesh Gupta's row needed a manual fix-up, so that cell is flagged
# with a distinct look: left aligned, plain black Arial text on a
# solid white background.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "B.Tech IT"

$d3 = $ws.Range("D3")
$d3.Font.Color = 0
$d3.Font.Name = "Arial"
$d3.Interior.Color = 16777215
$d3.HorizontalAlignment = -4131  # xlLeft
$d3.Value = "B.Tech CSE"

$ws.Range("D4").Value = "MCA"
